# Commit: "updated slides for better flow and fixing a few inconsistencies"
#
# Slide 10 ("Wrapper Class Example") has a bullet that originally read
# "Line #6 declares the a1 (an object instance) to be an ArrayList ...".
# The code sample actually puts that declaration on line #5, so the
# author corrected the text to say "Line #5" instead of "Line #6".
#
# We locate the exact run via the shape's TextFrame.TextRange, find the
# character span that holds the old sentence, and overwrite just that
# span so the surrounding runs/formatting are left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item(5)          # "Text Placeholder 2"
$tf = $sh.TextFrame
$tr = $tf.TextRange

$old = "Line #6 declares the a1 (an object instance) to be an "
$new = "Line #5 declares the a1 (an object instance) to be an "

$fullText = $tr.Text
$pos0 = $fullText.IndexOf($old)

if ($pos0 -ge 0) {
    $start = $pos0 + 1          # TextRange.Characters is 1-based
    $len = $old.Length
    $sub = $tr.Characters($start, $len)
    $sub.Text = $new
}
